# ZSS-1338 client test case - partial overflow test sheet rebuild
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Start clean: remove all existing content & the old merge ----
$ws.Cells.Clear()
$ws.Range("A1:B1").UnMerge()

# ---- Column widths / row heights ----
$ws.Columns.Item(1).ColumnWidth = 16.83
$ws.Columns.Item(7).ColumnWidth = 11.42
$ws.Rows.Item(6).RowHeight = 33

# ---- Cell values ----
$ws.Range("A1").Value = "normal"
$ws.Range("C1").Value = "a loooong text with left alingment"

$ws.Range("A2").Value = "no overflow, merged"
$ws.Range("C2").Value = "a loooong text with left alingment"

$ws.Range("A5").Value = "normal"
$ws.Range("F5").Value = "F5 loooong text with right alignment"

$ws.Range("A6").Value = "normal, higher"
$ws.Range("F6").Value = "F6 loooong text with right alignment"

$ws.Range("F7").Value = "F7 loooong text with right alignment"
$ws.Range("G7").Value = "right aligned"

$ws.Range("A8").Value = "cut next"
$ws.Range("E8").Value = "cut"
$ws.Range("F8").Value = "F8 loooong text with right alignment"

$ws.Range("A9").Value = "cut next 2"
$ws.Range("F9").Value = "F9 loooong text with right alignment"

$ws.Range("A10").Value = "no overflow"
$ws.Range("F10").Value = "right"

$ws.Range("A11").Value = "no overflow, merged"
$ws.Range("E11").Value = "F9 loooong text with right alignment"

# ---- Fonts: bold CJK font, applied over contiguous blocks to avoid
#      duplicate intermediate style records ----
$boldNormal1 = $ws.Range("A1:A2")
$boldNormal1.Font.Name = "新細明體"
$boldNormal1.Font.Bold = $true

$boldNormal2 = $ws.Range("A8:A11")
$boldNormal2.Font.Name = "新細明體"
$boldNormal2.Font.Bold = $true

# Bold CJK font, left aligned
$boldLeft = $ws.Range("A5:A7")
$boldLeft.Font.Name = "新細明體"
$boldLeft.Font.Bold = $true
$boldLeft.HorizontalAlignment = -4131

# ---- Alignment-only groups ----
$rightFill1 = $ws.Range("F5:F10")
$rightFill1.HorizontalAlignment = -4152

$ws.Range("G7").HorizontalAlignment = -4152

$leftGrp = $ws.Range("C2:D2")
$leftGrp.HorizontalAlignment = -4131

$rightGrp = $ws.Range("E11:F11")
$rightGrp.HorizontalAlignment = -4152

# ---- Merges ----
$ws.Range("C2:D2").Merge()
$ws.Range("E11:F11").Merge()

# ---- View settings ----
$excel.ActiveWindow.Zoom = 150
$ws.Range("D9").Select()
